$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New relationship rows being appended to the table (Id=14, 15), mirroring
# the existing "Id / ParentRelationshipId / ChildRelationshipId" rows above.
# Format the range as Text first so Excel stores these numeric-looking
# values ("14", "3", ...) as text (matching the rest of the sheet, which
# is entirely string-typed), then clear the formatting overrides so the
# new cells fall back to the sheet's default style, just like every other
# cell in the table.
$ws.Range("A15:C16").NumberFormat = "@"

$ws.Range("A15").Value = "14"
$ws.Range("B15").Value = "3"
$ws.Range("C15").Value = "3"

$ws.Range("A16").Value = "15"
$ws.Range("B16").Value = "3"
$ws.Range("C16").Value = "4"

$ws.Range("A15:C16").ClearFormats()
